$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), copying the formatting used by the
# other header cells (bold, bordered, centered) from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
